# Briefing Dra. Ana Maria - update metrics
# The "Dia 15/10" time estimate changes from 2hr to 3hr (1 dia).
$d = $word.ActiveDocument

$d.Content.Find.Execute("Dia 15/10: 2hr (1 dia)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Dia 15/10: 3hr (1 dia)", 2)
